$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.643.79"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.473.65"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "2.854.78"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").Value = "2.482.91"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").Value = "41.596.78"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0771"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").Value = "1.985.42"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0286"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("D48").Value = "2.710.70"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.52%  "
